# Update NNS values and delete econ uncertainty values
# (autumn/xls/data_default.xlsx, sheet "constants")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# --- Update NNS values (program_nns_xpertacf_smearpos / smearneg) ---
# Row 43: program_nns_xpertacf_smearpos  B: 1 -> 40, new comment in F
$ws.Range("B43").Value = 40
$ws.Range("F43").Value = "Emm's suggestion"

# Row 44: program_nns_xpertacf_smearneg  B: 2 -> 50, new comment in F
$ws.Range("B44").Value = 50
$ws.Range("F44").Value = "It needs to be higherd than nns_smearpos"

# --- Delete econ uncertainty (low/high) values, keep formatting ---
$ws.Range("C110:D110").ClearContents()
$ws.Range("C114:D114").ClearContents()
$ws.Range("C118:D118").ClearContents()
$ws.Range("C126:D126").ClearContents()
$ws.Range("C130:D130").ClearContents()
$ws.Range("C134:D134").ClearContents()

# --- Restore selection to match the author's last cursor position ---
$ws.Activate()
$ws.Range("C112").Select()
